$d = $word.ActiveDocument

# --- Change 1: restructure "a resident of Guiwan, this city" text with proofErr markers
#     and move the "Address"/"address" bookmark to be empty, right after "of ".
$bmAddr = $d.Bookmarks("Address")
$bmReq = $d.Bookmarks("requestor")
$rng = $d.Range($bmAddr.Start - 3, $bmReq.Start)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/></w:rPr><w:t xml:space="preserve">of </w:t></w:r><w:bookmarkStart w:id="1" w:name="address"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/></w:rPr><w:t>,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/></w:rPr><w:t>Guiwan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/></w:rPr><w:t xml:space="preserve">, this city, and the family is being indigent seeks Medical and Financial Assistance; this case is being referred to your office, </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)

# --- Change 2: swap bookmark names "Name2"/"requestee" -> "requestee"/"name2"
$bmName2 = $d.Bookmarks("Name2")
$name2Start = $bmName2.Start
$name2End = $bmName2.End
$name2Range = $d.Range($name2Start, $name2End)

$bmRequestee = $d.Bookmarks("requestee")
$requesteeStart = $bmRequestee.Start
$requesteeEnd = $bmRequestee.End
$requesteeRange = $d.Range($requesteeStart, $requesteeEnd)

$bmName2.Delete()
$bmRequestee.Delete()

$d.Bookmarks.Add("requestee", $requesteeRange)
$d.Bookmarks.Add("name2", $name2Range)
